# Refresh the cryptocurrency price / 1h-volume-change snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.882.56"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "2.455.96"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'517.07"
$ws.Range("E5").Value = "  -2.93%  "
$ws.Range("D6").Value = "'132.33"
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").Value = "2.460.88"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  -3.70%  "
$ws.Range("D11").Value = "'0.157"
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("D14").Value = "2.889.24"
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").Value = "57.799.38"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D16").Value = "'21.87"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("D18").Value = "2.458.28"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "'10.55"
$ws.Range("E19").Value = "  -4.00%  "
$ws.Range("D20").Value = "'318.65"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").Value = "'4.11"
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D23").Value = "'5.67"
$ws.Range("D24").Value = "'64.42"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").Value = "'0.406"
$ws.Range("E25").Value = "  -2.96%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0738"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'168.20"
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("D32").Value = "'6.20"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'17.93"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.33"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("D39").Value = "'36.30"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E40").Value = "  -4.54%  "
$ws.Range("D41").Value = "'0.784"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").Value = "'3.41"
$ws.Range("E42").Value = "  -4.12%  "
$ws.Range("D43").Value = "'270.23"
$ws.Range("E43").Value = "  -3.56%  "
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D46").Value = "'122.89"
$ws.Range("E46").Value = "  -4.94%  "
$ws.Range("D47").Value = "'0.0904"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").Value = "'0.0482"
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").Value = "1.722.07"
$ws.Range("E51").Value = "  -1.62%  "
